# "Báo cáo cá nhân" Lương sheet: remove the "Phụ cấp tại LONG XUYÊN" line
# (the template restructured LONG XUYÊN pay into "Tổng công" + "Lương cơ bản"
# only, dropping the separate phụ cấp line) and refresh the computed totals
# that depend on it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

# Delete row 13 ("Phụ cấp tại LONG XUYÊN" / 70000) entirely; everything below
# shifts up by one row and the sheet dimension shrinks from B35 to B34.
$ws.Rows.Item(13).Delete()

# "Tổng công tại LONG XUYÊN" (row 12, unaffected by the shift) goes from 2 -> 3.
$ws.Cells.Item(12, 2).Value = 3

# "Lương cơ bản tại LONG XUYÊN" (now row 13 after the shift up) is recomputed.
$ws.Cells.Item(13, 2).Value = 428571.4285714286

# "Tổng lương tại LONG XUYÊN" (now row 32) and "Tổng lương tại HỆ THỐNG"
# (now row 34) both pick up the same new total.
$ws.Cells.Item(32, 2).Value = 528571.4285714286
$ws.Cells.Item(34, 2).Value = 528571.4285714286
